# Fill in the test-case data (rows 6-10, columns B:I) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# input values for B:F per test-case row
$testData = @{
    6  = @(8, 12, 126, 333100360, 5)
    7  = @(50, 10, 1, 333100360, 5)
    8  = @(8, 12, 200, 333100360, 5)
    9  = @(8, 100, 1, 333100360, 5)
    10 = @(4, 12, 20, 333100360, 5)
}

foreach ($row in $testData.Keys) {
    $vals = $testData[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]

    $ws.Range("G$row").Formula = "=(`$B`$2/B$row + `$B`$2/D$row-`$B`$2/C$row) * F$row"
    $ws.Range("H$row").Formula = "=E$row+G$row"
    $ws.Range("I$row").Formula = "=IF(H$row > E$row, ""Increased"", ""Decreased"")"
}

# Restore the active selection to the single cell G9, as left by the author.
$null = $ws.Range("G9").Select()
